$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 385 (shifts old rows 385-452 down to 387-454)
$ws.Range("A385:A386").EntireRow.Insert()

# New row 385: Especial quality, Provincia de Melipilla record
$ws.Range("A385").Value = 4
$ws.Range("B385").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C385").Value = "Los Lagos"
$ws.Range("D385").Value = 45218
$ws.Range("E385").Value = 10
$ws.Range("F385").Value = "Fruta"
$ws.Range("G385").Value = 100101
$ws.Range("H385").Value = "Berries"
$ws.Range("I385").Value = 100112025
$ws.Range("J385").Value = "Frutilla"
$ws.Range("K385").Value = "Sin especificar"
$ws.Range("L385").Value = "Especial"
$ws.Range("M385").Value = 250
$ws.Range("N385").Value = 13000
$ws.Range("O385").Value = 13000
$ws.Range("P385").Value = 13000
$ws.Range("Q385").Value = "$/bandeja 7 kilos"
$ws.Range("R385").Value = "Provincia de Melipilla"
$ws.Range("S385").Value = 1857
$ws.Range("T385").Value = 7

# New row 386: Primera quality, Provincia de Melipilla record
$ws.Range("A386").Value = 4
$ws.Range("B386").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C386").Value = "Los Lagos"
$ws.Range("D386").Value = 45218
$ws.Range("E386").Value = 10
$ws.Range("F386").Value = "Fruta"
$ws.Range("G386").Value = 100101
$ws.Range("H386").Value = "Berries"
$ws.Range("I386").Value = 100112025
$ws.Range("J386").Value = "Frutilla"
$ws.Range("K386").Value = "Sin especificar"
$ws.Range("L386").Value = "Primera"
$ws.Range("M386").Value = 300
$ws.Range("N386").Value = 10000
$ws.Range("O386").Value = 10000
$ws.Range("P386").Value = 10000
$ws.Range("Q386").Value = "$/bandeja 7 kilos"
$ws.Range("R386").Value = "Provincia de Melipilla"
$ws.Range("S386").Value = 1429
$ws.Range("T386").Value = 7
